$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: the per-row Fecha / Volumen / Precio fields have been
# reshuffled across the existing data rows (row identity columns such as
# Mercado, Region, Categoria, etc. stay put). Apply the new values row by
# row for columns D, J, K, L, M, P.

$rows = @{
    2  = @(44413, 25, 24000, 25000, 24480, 1632)
    3  = @(44432, 34, 24000, 25000, 24500, 1633)
    4  = @(44428, 16, 25000, 26000, 25500, 1700)
    5  = @(44421, 18, 24000, 25000, 24500, 1633)
    6  = @(44349, 21, 24000, 25000, 24524, 1635)
    7  = @(44329, 25, 23000, 23000, 23000, 1533)
    8  = @(44435, 34, 24000, 25000, 24500, 1633)
    9  = @(44449, 18, 24000, 25000, 24500, 1633)
    10 = @(44446, 34, 24000, 25000, 24500, 1633)
    14 = @(44383, 25, 13000, 14000, 13480, 899)
    15 = @(44336, 34, 24000, 25000, 24500, 1633)
    17 = @(44400, 16, 24000, 25000, 24500, 1633)
    18 = @(44442, 28, 24000, 25000, 24500, 1633)
    19 = @(44453, 25, 25000, 26000, 25520, 1701)
    20 = @(44418, 16, 25000, 26000, 25500, 1700)
    21 = @(44343, 26, 23000, 24000, 23500, 1567)
    22 = @(44390, 34, 24000, 25000, 24500, 1633)
    23 = @(44397, 34, 23000, 24000, 23500, 1567)
    24 = @(44351, 34, 24000, 25000, 24500, 1633)
    25 = @(44411, 34, 25000, 26000, 25500, 1700)
    26 = @(44425, 25, 24000, 25000, 24520, 1635)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $vals[1]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals[2]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals[3]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals[5]   # P - Precio $/Kg
}
